$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range('D2:E2')
$rng.NumberFormat = "@"
$ws.Range('D2').Value = '255.06'
$ws.Range('E2').Value = '3.61%'
$rng.Style = "Normal"

$rng = $ws.Range('D3:E3')
$rng.NumberFormat = "@"
$ws.Range('D3').Value = '28.17'
$ws.Range('E3').Value = '-5.92%'
$rng.Style = "Normal"

$rng = $ws.Range('D4:E4')
$rng.NumberFormat = "@"
$ws.Range('D4').Value = '5.323'
$ws.Range('E4').Value = '3.05%'
$rng.Style = "Normal"

$rng = $ws.Range('D5:E5')
$rng.NumberFormat = "@"
$ws.Range('D5').Value = '0.05849'
$ws.Range('E5').Value = '0.88%'
$rng.Style = "Normal"

$rng = $ws.Range('D6:E6')
$rng.NumberFormat = "@"
$ws.Range('D6').Value = '6.707'
$ws.Range('E6').Value = '0.43%'
$rng.Style = "Normal"

$rng = $ws.Range('B7:E7')
$rng.NumberFormat = "@"
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '3.231'
$ws.Range('E7').Value = '0.48%'
$rng.Style = "Normal"

$rng = $ws.Range('B8:E8')
$rng.NumberFormat = "@"
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = '0.8675'
$ws.Range('E8').Value = '1.96%'
$rng.Style = "Normal"

$rng = $ws.Range('B9:E9')
$rng.NumberFormat = "@"
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').Value = '0.9094'
$ws.Range('E9').Value = '5.40%'
$rng.Style = "Normal"

$rng = $ws.Range('B10:E10')
$rng.NumberFormat = "@"
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').Value = '0.01065'
$ws.Range('E10').Value = '1,679.36%'
$rng.Style = "Normal"

$rng = $ws.Range('B11:E11')
$rng.NumberFormat = "@"
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1419'
$ws.Range('E11').Value = '2.99%'
$rng.Style = "Normal"

$rng = $ws.Range('B12:E12')
$rng.NumberFormat = "@"
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.07156'
$ws.Range('E12').Value = '0.17%'
$rng.Style = "Normal"

$rng = $ws.Range('B13:E13')
$rng.NumberFormat = "@"
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.03177'
$ws.Range('E13').Value = '-0.62%'
$rng.Style = "Normal"

$rng = $ws.Range('B14:E14')
$rng.NumberFormat = "@"
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.09224'
$ws.Range('E14').Value = '-1.61%'
$rng.Style = "Normal"

$rng = $ws.Range('B15:E15')
$rng.NumberFormat = "@"
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001547'
$ws.Range('E15').Value = '1.01%'
$rng.Style = "Normal"

$rng = $ws.Range('B16:E16')
$rng.NumberFormat = "@"
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005829'
$ws.Range('E16').Value = '-0.82%'
$rng.Style = "Normal"

$rng = $ws.Range('B17:E17')
$rng.NumberFormat = "@"
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.498'
$ws.Range('E17').Value = '-0.04%'
$rng.Style = "Normal"

$rng = $ws.Range('D18:E18')
$rng.NumberFormat = "@"
$ws.Range('D18').Value = '2.202'
$ws.Range('E18').Value = '-0.08%'
$rng.Style = "Normal"

$rng = $ws.Range('D19:E19')
$rng.NumberFormat = "@"
$ws.Range('D19').Value = '0.3173'
$ws.Range('E19').Value = '-0.73%'
$rng.Style = "Normal"

$rng = $ws.Range('E20')
$rng.NumberFormat = "@"
$ws.Range('E20').Value = '2.58%'
$rng.Style = "Normal"

$rng = $ws.Range('D21:E21')
$rng.NumberFormat = "@"
$ws.Range('D21').Value = '0.1315'
$ws.Range('E21').Value = '1.18%'
$rng.Style = "Normal"

$rng = $ws.Range('D22:E22')
$rng.NumberFormat = "@"
$ws.Range('D22').Value = '3.534'
$ws.Range('E22').Value = '1.04%'
$rng.Style = "Normal"

$rng = $ws.Range('D23:E23')
$rng.NumberFormat = "@"
$ws.Range('D23').Value = '0.04150'
$ws.Range('E23').Value = '0.11%'
$rng.Style = "Normal"

$rng = $ws.Range('E24')
$rng.NumberFormat = "@"
$ws.Range('E24').Value = '-0.12%'
$rng.Style = "Normal"

$rng = $ws.Range('D25')
$rng.NumberFormat = "@"
$ws.Range('D25').Value = '0.005041'
$rng.Style = "Normal"

$rng = $ws.Range('D26:E26')
$rng.NumberFormat = "@"
$ws.Range('D26').Value = '0.001226'
$ws.Range('E26').Value = '-0.06%'
$rng.Style = "Normal"

$rng = $ws.Range('D27:E27')
$rng.NumberFormat = "@"
$ws.Range('D27').Value = '0.0001699'
$ws.Range('E27').Value = '41.63%'
$rng.Style = "Normal"

$rng = $ws.Range('E28')
$rng.NumberFormat = "@"
$ws.Range('E28').Value = '33.73%'
$rng.Style = "Normal"

$rng = $ws.Range('D40:E40')
$rng.NumberFormat = "@"
$ws.Range('D40').Value = '0.03855'
$ws.Range('E40').Value = '2.50%'
$rng.Style = "Normal"

$rng = $ws.Range('D41:E41')
$rng.NumberFormat = "@"
$ws.Range('D41').Value = '0.1100'
$ws.Range('E41').Value = '2.80%'
$rng.Style = "Normal"

$rng = $ws.Range('D42:E42')
$rng.NumberFormat = "@"
$ws.Range('D42').Value = '0.002199'
$ws.Range('E42').Value = '-0.02%'
$rng.Style = "Normal"

$rng = $ws.Range('D43:E43')
$rng.NumberFormat = "@"
$ws.Range('D43').Value = '0.002948'
$ws.Range('E43').Value = '-48.53%'
$rng.Style = "Normal"

$rng = $ws.Range('D44:E44')
$rng.NumberFormat = "@"
$ws.Range('D44').Value = '0.01099'
$ws.Range('E44').Value = '15.02%'
$rng.Style = "Normal"

$rng = $ws.Range('D45')
$rng.NumberFormat = "@"
$ws.Range('D45').Value = '0.00005226'
$rng.Style = "Normal"

$rng = $ws.Range('D47:E47')
$rng.NumberFormat = "@"
$ws.Range('D47').Value = '0.08750'
$ws.Range('E47').Value = '50.91%'
$rng.Style = "Normal"

$rng = $ws.Range('D48:E48')
$rng.NumberFormat = "@"
$ws.Range('D48').Value = '0.002156'
$ws.Range('E48').Value = '-1.19%'
$rng.Style = "Normal"
